$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 202; this shifts existing rows 202:244 down to 203:245
$ws.Rows(202).Insert()

# Populate the newly inserted row 202 with the new weekly record
$ws.Cells.Item(202, 1).Value = 5
$ws.Cells.Item(202, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(202, 3).Value = "Maule"
$ws.Cells.Item(202, 4).Value = "2022-05-13"
$ws.Cells.Item(202, 5).Value = 7
$ws.Cells.Item(202, 6).Value = 100112009
$ws.Cells.Item(202, 7).Value = "Acelga"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 500
$ws.Cells.Item(202, 11).Value = 3000
$ws.Cells.Item(202, 12).Value = 3000
$ws.Cells.Item(202, 13).Value = 3000
$ws.Cells.Item(202, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(202, 15).Value = "Región del Maule"
$ws.Cells.Item(202, 16).Value = 750
$ws.Cells.Item(202, 17).Value = 4
$ws.Cells.Item(202, 18).Value = "Hortaliza"
